$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 90.59999999999999
$ws.Range("I9").Value = 65.57143000000001
$ws.Range("J9").Value = 149
$ws.Range("K9").Value = 65.57143000000001
$ws.Range("L9").Value = 149
$ws.Range("M9").Value = 103.42857
$ws.Range("N9").Value = -487
$ws.Range("H19").Value = 1603.7778
$ws.Range("I19").Value = 244.125
$ws.Range("K19").Value = 244.125
$ws.Range("M19").Value = -69.125
$ws.Range("H40").Value = 2427.1428
$ws.Range("I40").Value = 1997.75
$ws.Range("J40").Value = 2999.6667
$ws.Range("K40").Value = 1997.75
$ws.Range("L40").Value = 2999.6667
$ws.Range("M40").Value = -1822.75
$ws.Range("N40").Value = -3349.6667
$ws.Range("H82").Value = 495.5
$ws.Range("I82").Value = 495.5
$ws.Range("K82").Value = 1486.5
$ws.Range("M82").Value = -1080.5
$ws.Range("H85").Value = 495.5
$ws.Range("I85").Value = 495.5
$ws.Range("K85").Value = 1486.5
$ws.Range("M85").Value = -82.5
$ws.Range("H135").Value = 1678.5385
$ws.Range("I135").Value = 1165.6
$ws.Range("J135").Value = 3388.3333
$ws.Range("K135").Value = 10490.4
$ws.Range("L135").Value = 30494.9997
$ws.Range("M135").Value = -7955.4
$ws.Range("N135").Value = -35564.9997

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = $null
$ws.Range("N17").Value = 0
$ws.Range("H32").Value = 5203.5
$ws.Range("I32").Value = 3421.25
$ws.Range("K32").Value = 3421.25
$ws.Range("M32").Value = -3134.25
$ws.Range("H63").Value = 1933.3334
$ws.Range("I63").Value = 1933.3334
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1933.3334
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = -1247.3334
$ws.Range("H66").Value = 1933.3334
$ws.Range("I66").Value = 1933.3334
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9666.666999999999
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = -6234.666999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1217
$ws.Range("I20").Value = 1285.5555
$ws.Range("J20").Value = 908.5
$ws.Range("K20").Value = 1285.5555
$ws.Range("L20").Value = 908.5
$ws.Range("M20").Value = -1038.5555
$ws.Range("N20").Value = -1402.5
$ws.Range("H86").Value = 1070.9
$ws.Range("I86").Value = 868.1667
$ws.Range("K86").Value = 868.1667
$ws.Range("M86").Value = 254.8333
$ws.Range("H89").Value = 1070.9
$ws.Range("I89").Value = 868.1667
$ws.Range("K89").Value = 4340.8335
$ws.Range("M89").Value = 1275.1665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 89293.336
$ws.Range("I22").Value = 132877.67
$ws.Range("J22").Value = 23916.834
$ws.Range("K22").Value = 132877.67
$ws.Range("L22").Value = 23916.834
$ws.Range("M22").Value = -132527.67
$ws.Range("N22").Value = -24616.834
$ws.Range("H99").Value = 15044.333
$ws.Range("I99").Value = 9509.700000000001
$ws.Range("J99").Value = 18997.643
$ws.Range("K99").Value = 9509.700000000001
$ws.Range("L99").Value = 18997.643
$ws.Range("M99").Value = -8011.700000000001
$ws.Range("N99").Value = -21993.643
$ws.Range("H126").Value = 15044.333
$ws.Range("I126").Value = 9509.700000000001
$ws.Range("J126").Value = 18997.643
$ws.Range("K126").Value = 28529.1
$ws.Range("L126").Value = 56992.929
$ws.Range("M126").Value = -26059.1
$ws.Range("N126").Value = -61932.929
$ws.Range("H134").Value = 2936
$ws.Range("I134").Value = 2493.9167
$ws.Range("K134").Value = 7481.750100000001
$ws.Range("M134").Value = -4946.750100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 412
$ws.Range("I5").Value = 433
$ws.Range("J5").Value = 349
$ws.Range("K5").Value = 1299
$ws.Range("L5").Value = 1047
$ws.Range("M5").Value = -1187
$ws.Range("N5").Value = -1271
$ws.Range("H23").Value = 252525.5
$ws.Range("I23").Value = 101
$ws.Range("J23").Value = 336667
$ws.Range("K23").Value = 303
$ws.Range("L23").Value = 1010001
$ws.Range("M23").Value = -68
$ws.Range("N23").Value = -1010471
$ws.Range("H34").Value = 2523
$ws.Range("J34").Value = 3750
$ws.Range("L34").Value = 11250
$ws.Range("N34").Value = -11418
$ws.Range("H39").Value = 10631.23
$ws.Range("J39").Value = 10631.23
$ws.Range("L39").Value = 31893.69
$ws.Range("N39").Value = -32481.69
$ws.Range("H40").Value = 44.666668
$ws.Range("I40").Value = 26.555555
$ws.Range("K40").Value = 106.22222
$ws.Range("M40").Value = -37.22221999999999
$ws.Range("H51").Value = 933
$ws.Range("I51").Value = 1199
$ws.Range("J51").Value = 800
$ws.Range("K51").Value = 3597
$ws.Range("L51").Value = 2400
$ws.Range("M51").Value = -3137
$ws.Range("N51").Value = -3320
$ws.Range("H55").Value = 7518.625
$ws.Range("I55").Value = 3124.5
$ws.Range("J55").Value = 8983.333000000001
$ws.Range("K55").Value = 9373.5
$ws.Range("L55").Value = 26949.999
$ws.Range("M55").Value = -9196.5
$ws.Range("N55").Value = -27303.999
$ws.Range("H94").Value = 3046.5715
$ws.Range("J94").Value = 4105.2
$ws.Range("L94").Value = 12315.6
$ws.Range("N94").Value = -13667.6
$ws.Range("H135").Value = 412
$ws.Range("I135").Value = 433
$ws.Range("J135").Value = 349
$ws.Range("K135").Value = 3897
$ws.Range("L135").Value = 3141
$ws.Range("M135").Value = -1362
$ws.Range("N135").Value = -8211

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 1388800.8
$ws.Range("J18").Value = 18400.166
$ws.Range("L18").Value = 18400.166
$ws.Range("N18").Value = -18986.166
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = $null
$ws.Range("N34").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = $null
$ws.Range("N79").Value = 0
$ws.Range("H80").Value = 3797.4285
$ws.Range("I80").Value = 3597
$ws.Range("K80").Value = 3597
$ws.Range("M80").Value = -2599
$ws.Range("H83").Value = 3797.4285
$ws.Range("I83").Value = 3597
$ws.Range("K83").Value = 17985
$ws.Range("M83").Value = -12993

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 20026600
$ws.Range("I23").Value = 25019500
$ws.Range("J23").Value = 55000
$ws.Range("K23").Value = 25019500
$ws.Range("L23").Value = 55000
$ws.Range("M23").Value = -25019270
$ws.Range("N23").Value = -55460
$ws.Range("H46").Value = 3993.6667
$ws.Range("I46").Value = 2999.923
$ws.Range("J46").Value = 6577.4
$ws.Range("K46").Value = 2999.923
$ws.Range("L46").Value = 6577.4
$ws.Range("M46").Value = -2811.923
$ws.Range("N46").Value = -6953.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4200300
$ws.Range("H100").Value = 1406.2142
$ws.Range("I100").Value = 1393.4445
$ws.Range("J100").Value = 1429.2
$ws.Range("K100").Value = 2786.889
$ws.Range("L100").Value = 2858.4
$ws.Range("M100").Value = -2245.889
$ws.Range("N100").Value = -3940.4
